$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 934
$ws.Range("I8").Value = 107.09091
$ws.Range("J8").Value = 2450
$ws.Range("K8").Value = 321.27273
$ws.Range("L8").Value = 7350
$ws.Range("M8").Value = -182.27273
$ws.Range("N8").Value = -7628
$ws.Range("H9").Value = 4762205.5
$ws.Range("I9").Value = 10000116
$ws.Range("J9").Value = 468.54544
$ws.Range("K9").Value = 10000116
$ws.Range("L9").Value = 468.54544
$ws.Range("M9").Value = -9999947
$ws.Range("N9").Value = -806.54544
$ws.Range("H97").Value = 2483.25
$ws.Range("J97").Value = 2483.25
$ws.Range("L97").Value = 7449.75
$ws.Range("N97").Value = -8441.75
$ws.Range("H112").Value = 4808910
$ws.Range("I112").Value = 3372.5
$ws.Range("J112").Value = 5209371.5
$ws.Range("K112").Value = 10117.5
$ws.Range("L112").Value = 15628114.5
$ws.Range("M112").Value = -9009.5
$ws.Range("N112").Value = -15630330.5
$ws.Range("H132").Value = 4083176
$ws.Range("I132").Value = 4762938.5
$ws.Range("J132").Value = 4599.4287
$ws.Range("K132").Value = 14288815.5
$ws.Range("L132").Value = 13798.2861
$ws.Range("M132").Value = -14286285.5
$ws.Range("N132").Value = -18858.2861

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3939.5
$ws.Range("I32").Value = 2629.4492
$ws.Range("K32").Value = 2629.4492
$ws.Range("M32").Value = -2342.4492
$ws.Range("H45").Value = 2123.524
$ws.Range("I45").Value = 1040
$ws.Range("K45").Value = 1040
$ws.Range("M45").Value = -663
$ws.Range("H61").Value = 3178.9355
$ws.Range("I61").Value = 1496.8823
$ws.Range("J61").Value = 5221.4287
$ws.Range("K61").Value = 1496.8823
$ws.Range("L61").Value = 5221.4287
$ws.Range("M61").Value = -1284.8823
$ws.Range("N61").Value = -5645.4287
$ws.Range("H74").Value = 826.95654
$ws.Range("I74").Value = 791
$ws.Range("K74").Value = 791
$ws.Range("M74").Value = 83
$ws.Range("H77").Value = 826.95654
$ws.Range("I77").Value = 791
$ws.Range("K77").Value = 3955
$ws.Range("M77").Value = 413
$ws.Range("H122").Value = 2671.25
$ws.Range("I122").Value = 1687.1765
$ws.Range("J122").Value = 5061.143
$ws.Range("K122").Value = 5061.529500000001
$ws.Range("L122").Value = 15183.429
$ws.Range("M122").Value = -2611.529500000001
$ws.Range("N122").Value = -20083.429
$ws.Range("H125").Value = 29642.857
$ws.Range("J125").Value = 29642.857
$ws.Range("L125").Value = 29642.857
$ws.Range("N125").Value = -39482.857
$ws.Range("H132").Value = 13701480
$ws.Range("I132").Value = 17244052
$ws.Range("J132").Value = 3533.2
$ws.Range("K132").Value = 51732156
$ws.Range("L132").Value = 10599.6
$ws.Range("M132").Value = -51729626
$ws.Range("N132").Value = -15659.6
$ws.Range("H136").Value = 3178.9355
$ws.Range("I136").Value = 1496.8823
$ws.Range("J136").Value = 5221.4287
$ws.Range("K136").Value = 4490.6469
$ws.Range("L136").Value = 15664.2861
$ws.Range("M136").Value = -1940.6469
$ws.Range("N136").Value = -20764.2861
$ws.Range("H139").Value = 24750
$ws.Range("J139").Value = 24750
$ws.Range("L139").Value = 24750
$ws.Range("N139").Value = -35030

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1568.3871
$ws.Range("I99").Value = 1161.2609
$ws.Range("J99").Value = 2738.875
$ws.Range("K99").Value = 1161.2609
$ws.Range("L99").Value = 2738.875
$ws.Range("M99").Value = 336.7391
$ws.Range("N99").Value = -5734.875
$ws.Range("H105").Value = 1616.5
$ws.Range("I105").Value = 1223.7778
$ws.Range("J105").Value = 2058.3125
$ws.Range("K105").Value = 1223.7778
$ws.Range("L105").Value = 2058.3125
$ws.Range("M105").Value = 523.2221999999999
$ws.Range("N105").Value = -5552.3125
$ws.Range("H134").Value = 2185.0244
$ws.Range("I134").Value = 1658
$ws.Range("J134").Value = 4744.857
$ws.Range("K134").Value = 4974
$ws.Range("L134").Value = 14234.571
$ws.Range("M134").Value = -2439
$ws.Range("N134").Value = -19304.571

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3536.5667
$ws.Range("I86").Value = 2531.6875
$ws.Range("J86").Value = 4685
$ws.Range("K86").Value = 2531.6875
$ws.Range("L86").Value = 4685
$ws.Range("M86").Value = -1408.6875
$ws.Range("N86").Value = -6931
$ws.Range("H89").Value = 3536.5667
$ws.Range("I89").Value = 2531.6875
$ws.Range("J89").Value = 4685
$ws.Range("K89").Value = 12658.4375
$ws.Range("L89").Value = 23425
$ws.Range("M89").Value = -7042.4375
$ws.Range("N89").Value = -34657
$ws.Range("H99").Value = 3703
$ws.Range("I99").Value = 1212
$ws.Range("K99").Value = 1212
$ws.Range("M99").Value = 286
$ws.Range("H124").Value = 28326
$ws.Range("J124").Value = 28326
$ws.Range("L124").Value = 28326
$ws.Range("N124").Value = -33236
$ws.Range("H126").Value = 3703
$ws.Range("I126").Value = 1212
$ws.Range("K126").Value = 3636
$ws.Range("M126").Value = -1166
$ws.Range("H132").Value = 2832.0222
$ws.Range("I132").Value = 1752.4231
$ws.Range("J132").Value = 4309.3687
$ws.Range("K132").Value = 5257.2693
$ws.Range("L132").Value = 12928.1061
$ws.Range("M132").Value = -2727.2693
$ws.Range("N132").Value = -17988.1061
$ws.Range("H134").Value = 1619.1316
$ws.Range("I134").Value = 1222.5667
$ws.Range("J134").Value = 3106.25
$ws.Range("K134").Value = 3667.7001
$ws.Range("L134").Value = 9318.75
$ws.Range("M134").Value = -1132.7001
$ws.Range("N134").Value = -14388.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1657.5698
$ws.Range("I68").Value = 660.5484
$ws.Range("J68").Value = 2219.5273
$ws.Range("K68").Value = 1981.6452
$ws.Range("L68").Value = 6658.581900000001
$ws.Range("M68").Value = -1170.6452
$ws.Range("N68").Value = -8280.581900000001
$ws.Range("H71").Value = 1657.5698
$ws.Range("I71").Value = 660.5484
$ws.Range("J71").Value = 2219.5273
$ws.Range("K71").Value = 5944.9356
$ws.Range("L71").Value = 19975.7457
$ws.Range("M71").Value = -1888.9356
$ws.Range("N71").Value = -28087.7457
$ws.Range("H94").Value = 3507.6667
$ws.Range("I94").Value = 1420
$ws.Range("J94").Value = 3925.2
$ws.Range("K94").Value = 4260
$ws.Range("L94").Value = 11775.6
$ws.Range("M94").Value = -3584
$ws.Range("N94").Value = -13127.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3534.7334
$ws.Range("I80").Value = 3172.2
$ws.Range("J80").Value = 4259.8
$ws.Range("K80").Value = 3172.2
$ws.Range("L80").Value = 4259.8
$ws.Range("M80").Value = -2174.2
$ws.Range("N80").Value = -6255.8
$ws.Range("H83").Value = 3534.7334
$ws.Range("I83").Value = 3172.2
$ws.Range("J83").Value = 4259.8
$ws.Range("K83").Value = 15861
$ws.Range("L83").Value = 21299
$ws.Range("M83").Value = -10869
$ws.Range("N83").Value = -31283
$ws.Range("H102").Value = 1824.2413
$ws.Range("I102").Value = 1162.579
$ws.Range("J102").Value = 3081.4
$ws.Range("K102").Value = 1162.579
$ws.Range("L102").Value = 3081.4
$ws.Range("M102").Value = 459.421
$ws.Range("N102").Value = -6325.4
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 3674.1292
$ws.Range("I132").Value = 2613.6191
$ws.Range("J132").Value = 5901.2
$ws.Range("K132").Value = 7840.8573
$ws.Range("L132").Value = 17703.6
$ws.Range("M132").Value = -5310.8573
$ws.Range("N132").Value = -22763.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2436
$ws.Range("I7").Value = 1658.8572
$ws.Range("K7").Value = 1658.8572
$ws.Range("M7").Value = -1546.8572
$ws.Range("H40").Value = 2542
$ws.Range("I40").Value = 1933
$ws.Range("K40").Value = 1933
$ws.Range("M40").Value = -1797
$ws.Range("H46").Value = 2293.3333
$ws.Range("J46").Value = 2996.6667
$ws.Range("L46").Value = 2996.6667
$ws.Range("N46").Value = -3372.6667
$ws.Range("H62").Value = 24900
$ws.Range("J62").Value = 24900
$ws.Range("L62").Value = 24900
$ws.Range("N62").Value = -26148
$ws.Range("H65").Value = 24900
$ws.Range("J65").Value = 24900
$ws.Range("L65").Value = 74700
$ws.Range("N65").Value = -80940
$ws.Range("H126").Value = 2436
$ws.Range("I126").Value = 1658.8572
$ws.Range("K126").Value = 4976.571599999999
$ws.Range("M126").Value = -2506.571599999999
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140
$ws.Range("H136").Value = 2001637
$ws.Range("I136").Value = 2632869.8
$ws.Range("J136").Value = 2733.75
$ws.Range("K136").Value = 7898609.399999999
$ws.Range("L136").Value = 8201.25
$ws.Range("M136").Value = -7896059.399999999
$ws.Range("N136").Value = -13301.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5380
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 6475
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 19425
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -23265
$ws.Range("H122").Value = 402309.28
$ws.Range("I122").Value = 556990.7
$ws.Range("J122").Value = 4557.143
$ws.Range("K122").Value = 1670972.1
$ws.Range("L122").Value = 13671.429
$ws.Range("M122").Value = -1668522.1
$ws.Range("N122").Value = -18571.429
$ws.Range("H126").Value = 2779695
$ws.Range("I126").Value = 1518.909
$ws.Range("K126").Value = 4556.727000000001
$ws.Range("M126").Value = -2086.727000000001
$ws.Range("H132").Value = 204054.42
$ws.Range("I132").Value = 271708.72
$ws.Range("K132").Value = 815126.1599999999
$ws.Range("M132").Value = -812596.1599999999
$ws.Range("H136").Value = 1046.4237
$ws.Range("I136").Value = 546.01886
$ws.Range("K136").Value = 1638.05658
$ws.Range("M136").Value = 911.9434200000001
